# "new experiments, data, analysis, camera ready"
#
# Adds:
#  - a new results column (I2:I9) next to the first "ant" condition block
#  - a third condition block (rows 21-29) for a new "buffao" animal, with
#    header row 21 (small / strong / wild) and 8 binary data rows
#
# Cell-write order below mirrors how the sheet was actually authored: the
# "buffao" label (row 22) was typed before the new "wild" header (row 21),
# which is why "buffao" lands at shared-string index 5 and "wild" at 6.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New column I: extra measurement for the first ("ant") block ---
$ws.Range("I2").Value = 0.127188041
$ws.Range("I3").Value = 0.085787136
$ws.Range("I4").Value = 0.069772448
$ws.Range("I5").Value = 0.060471832
$ws.Range("I6").Value = 0.274718337
$ws.Range("I7").Value = 0.155711647
$ws.Range("I8").Value = 0.122877287
$ws.Range("I9").Value = 0.103473271

# --- New third block: "buffao" / "wild" condition (rows 21-29) ---
$ws.Range("A22").Value = "buffao"

$ws.Range("C21").Value = "small"
$ws.Range("D21").Value = "strong"
$ws.Range("E21").Value = "wild"

$ws.Range("C22").Value = 0
$ws.Range("D22").Value = 1
$ws.Range("E22").Value = 1

$ws.Range("C23").Value = 0
$ws.Range("D23").Value = 1
$ws.Range("E23").Value = 0

$ws.Range("C24").Value = 0
$ws.Range("D24").Value = 0
$ws.Range("E24").Value = 1

$ws.Range("C25").Value = 0
$ws.Range("D25").Value = 0
$ws.Range("E25").Value = 0

$ws.Range("C26").Value = 1
$ws.Range("D26").Value = 1
$ws.Range("E26").Value = 1

$ws.Range("C27").Value = 1
$ws.Range("D27").Value = 1
$ws.Range("E27").Value = 0

$ws.Range("C28").Value = 1
$ws.Range("D28").Value = 0
$ws.Range("E28").Value = 1

$ws.Range("C29").Value = 1
$ws.Range("D29").Value = 0
$ws.Range("E29").Value = 0

# Match the author's final selection, left on the newly-added column.
$ws.Range("I2:I9").Select() | Out-Null
